# "Add also at to by institute author list"
#
# Five authors in the Imperial collaboration list gained a second
# ("also at") institutional affiliation, so their "Number of
# affiliations" count (column I) goes from 1 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imperial")
$ws.Select()

$rows = @(12, 14, 16, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = 2
}

# Reflect the view state left in the saved workbook: scrolled so
# column I is visible, with I21 as the active selection.
$ws.Range("I21").Select()
$excel.ActiveWindow.ScrollColumn = 9
